$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are forced to text via a leading apostrophe,
# mirroring Excel's own quote-prefix text-entry convention so that
# numeric-looking prices (e.g. "595.43") are not auto-converted to numbers.

$ws.Range('D2').Value = "'68.040.16"
$ws.Range('E2').Value = "  -0.32%  "
$ws.Range('D3').Value = "'3.659.75"
$ws.Range('E3').Value = "  -0.96%  "
$ws.Range('E4').Value = "  +0.17%  "
$ws.Range('D5').Value = "'595.43"
$ws.Range('E5').Value = "  +2.17%  "
$ws.Range('D6').Value = "'192.40"
$ws.Range('E6').Value = "  +4.46%  "
$ws.Range('D7').Value = "'0.620"
$ws.Range('E7').Value = "  -1.34%  "
$ws.Range('E8').Value = "  +0.23%  "
$ws.Range('D9').Value = "'0.696"
$ws.Range('E9').Value = "  -3.26%  "
$ws.Range('B10').Value = "Avalanche"
$ws.Range('C10').Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('D10').Value = "'56.79"
$ws.Range('E10').Value = "  +0.34%  "
$ws.Range('B11').Value = "Dogecoin"
$ws.Range('C11').Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range('D11').Value = "'0.151"
$ws.Range('E11').Value = "  -7.18%  "
$ws.Range('D12').Value = "'0.0000270"
$ws.Range('E12').Value = "  -7.07%  "
$ws.Range('D13').Value = "'10.19"
$ws.Range('E13').Value = "  -2.22%  "
$ws.Range('D14').Value = "'4.241.45"
$ws.Range('E14').Value = "  -0.93%  "
$ws.Range('D15').Value = "'3.664.33"
$ws.Range('E15').Value = "  -0.92%  "
$ws.Range('E16').Value = "  +0.33%  "
$ws.Range('D17').Value = "'18.85"
$ws.Range('E17').Value = "  -2.75%  "
$ws.Range('E18').Value = "  -1.79%  "
$ws.Range('D19').Value = "'67.910.20"
$ws.Range('E19').Value = "  -0.19%  "
$ws.Range('E20').Value = "  -2.89%  "
$ws.Range('D21').Value = "'399.61"
$ws.Range('E21').Value = "  -2.51%  "
$ws.Range('E22').Value = "  -1.77%  "
$ws.Range('D23').Value = "'87.57"
$ws.Range('E23').Value = "  -1.11%  "
$ws.Range('B24').Value = "RenderToken"
$ws.Range('C24').Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('D24').Value = "'11.09"
$ws.Range('E24').Value = "  +0.30%  "
$ws.Range('B25').Value = "ImmutableX"
$ws.Range('C25').Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('D25').Value = "'2.95"
$ws.Range('E25').Value = "  -2.80%  "
$ws.Range('D26').Value = "'12.54"
$ws.Range('E26').Value = "  -2.04%  "
$ws.Range('E27').Value = "  +0.25%  "
$ws.Range('E28').Value = "  -5.22%  "
$ws.Range('D29').Value = "'9.33"
$ws.Range('E29').Value = "  -1.68%  "
$ws.Range('D30').Value = "'31.80"
$ws.Range('E30').Value = "  -2.73%  "
$ws.Range('D31').Value = "'7.23"
$ws.Range('E31').Value = "  -2.63%  "
$ws.Range('E32').Value = "  -1.86%  "
$ws.Range('D33').Value = "'44.63"
$ws.Range('E33').Value = "  +1.50%  "
$ws.Range('D34').Value = "'66.14"
$ws.Range('E34').Value = "  +1.97%  "
$ws.Range('E35').Value = "  -1.19%  "
$ws.Range('D36').Value = "'608.15"
$ws.Range('E36').Value = "  +2.77%  "
$ws.Range('E37').Value = "  +0.11%  "
$ws.Range('E38').Value = "  -2.17%  "
$ws.Range('D39').Value = "'0.999"
$ws.Range('E39').Value = "  -0.23%  "
$ws.Range('D40').Value = "'0.0₃0766"
$ws.Range('E40').Value = "  -13.90%  "
$ws.Range('E41').Value = "  -0.94%  "
$ws.Range('D42').Value = "'2.88"
$ws.Range('E42').Value = "  -4.08%  "
$ws.Range('E43').Value = "  -2.16%  "
$ws.Range('E44').Value = "  -7.62%  "
$ws.Range('E45').Value = "  +0.84%  "
$ws.Range('D46').Value = "'2.770.08"
$ws.Range('E46').Value = "  -0.53%  "
$ws.Range('D47').Value = "'3.12"
$ws.Range('E47').Value = "  -0.23%  "
$ws.Range('D48').Value = "'143.30"
$ws.Range('E48').Value = "  +2.73%  "
$ws.Range('D49').Value = "'8.74"
$ws.Range('E49').Value = "  -6.13%  "
$ws.Range('E50').Value = "  -2.78%  "
$ws.Range('D51').Value = "'2.45"
$ws.Range('E51').Value = "  -15.73%  "
